# Auto-generated script to apply numeric updates to Cactuar_Profits workbook
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 239.75
$ws.Range("I12").Value = 239.75
$ws.Range("K12").Value = 239.75
$ws.Range("M12").Value = -69.75
$ws.Range("H18").Value = 166666850
$ws.Range("I18").Value = 166666850
$ws.Range("K18").Value = 166666850
$ws.Range("M18").Value = -166666566
$ws.Range("H62").Value = 2600.9
$ws.Range("I62").Value = 2335.6667
$ws.Range("J62").Value = 2998.75
$ws.Range("K62").Value = 2335.6667
$ws.Range("L62").Value = 2998.75
$ws.Range("M62").Value = -1711.6667
$ws.Range("N62").Value = -4246.75
$ws.Range("H65").Value = 2600.9
$ws.Range("I65").Value = 2335.6667
$ws.Range("J65").Value = 2998.75
$ws.Range("K65").Value = 11678.3335
$ws.Range("L65").Value = 14993.75
$ws.Range("M65").Value = -8558.333500000001
$ws.Range("N65").Value = -21233.75
$ws.Range("H76").Value = 3683.6924
$ws.Range("I76").Value = 3498.25
$ws.Range("K76").Value = 3498.25
$ws.Range("M76").Value = -3183.25
$ws.Range("H79").Value = 3683.6924
$ws.Range("I79").Value = 3498.25
$ws.Range("K79").Value = 3498.25
$ws.Range("M79").Value = -2406.25
$ws.Range("H103").Value = 1448.8334
$ws.Range("I103").Value = 1438.6
$ws.Range("J103").Value = 1500
$ws.Range("K103").Value = 4315.799999999999
$ws.Range("L103").Value = 4500
$ws.Range("M103").Value = -3729.799999999999
$ws.Range("N103").Value = -5672
$ws.Range("H111").Value = 1794.7142
$ws.Range("J111").Value = 763
$ws.Range("L111").Value = 2289
$ws.Range("N111").Value = -8423
$ws.Range("H137").Value = 3478.487
$ws.Range("I137").Value = 1192.05
$ws.Range("J137").Value = 5885.263
$ws.Range("K137").Value = 3576.15
$ws.Range("L137").Value = 17655.789
$ws.Range("M137").Value = -1026.15
$ws.Range("N137").Value = -22755.789
$ws.Range("H138").Value = 3566.577
$ws.Range("I138").Value = 1459
$ws.Range("J138").Value = 4682.353
$ws.Range("K138").Value = 4377
$ws.Range("L138").Value = 14047.059
$ws.Range("M138").Value = 763
$ws.Range("N138").Value = -24327.059
$ws.Range("H141").Value = 6204.9585
$ws.Range("I141").Value = 5405.5454
$ws.Range("K141").Value = 16216.6362
$ws.Range("M141").Value = -11036.6362

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4583.9106
$ws.Range("I32").Value = 2869.8293
$ws.Range("J32").Value = 9269.066000000001
$ws.Range("K32").Value = 2869.8293
$ws.Range("L32").Value = 9269.066000000001
$ws.Range("M32").Value = -2582.8293
$ws.Range("N32").Value = -9843.066000000001
$ws.Range("H61").Value = 10343.625
$ws.Range("I61").Value = 2305.2856
$ws.Range("K61").Value = 2305.2856
$ws.Range("M61").Value = -2093.2856
$ws.Range("H74").Value = 1198.6522
$ws.Range("I74").Value = 937.46155
$ws.Range("K74").Value = 937.46155
$ws.Range("M74").Value = -63.46154999999999
$ws.Range("H77").Value = 1198.6522
$ws.Range("I77").Value = 937.46155
$ws.Range("K77").Value = 4687.30775
$ws.Range("M77").Value = -319.3077499999999
$ws.Range("H132").Value = 22677.5
$ws.Range("I132").Value = 26282.088
$ws.Range("K132").Value = 78846.264
$ws.Range("M132").Value = -76316.264
$ws.Range("H136").Value = 10343.625
$ws.Range("I136").Value = 2305.2856
$ws.Range("K136").Value = 6915.8568
$ws.Range("M136").Value = -4365.8568

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 824.25
$ws.Range("I22").Value = 824.25
$ws.Range("K22").Value = 824.25
$ws.Range("M22").Value = -651.25
$ws.Range("H105").Value = 3053.375
$ws.Range("I105").Value = 1755
$ws.Range("J105").Value = 3238.8572
$ws.Range("K105").Value = 1755
$ws.Range("L105").Value = 3238.8572
$ws.Range("M105").Value = -8
$ws.Range("N105").Value = -6732.8572
$ws.Range("H107").Value = 1230.75
$ws.Range("I107").Value = 1211.5
$ws.Range("K107").Value = 1211.5
$ws.Range("M107").Value = 708.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2230.4
$ws.Range("J11").Value = 2762.5
$ws.Range("L11").Value = 2762.5
$ws.Range("N11").Value = -3042.5
$ws.Range("H22").Value = 975.8
$ws.Range("I22").Value = 714.5
$ws.Range("K22").Value = 714.5
$ws.Range("M22").Value = -364.5
$ws.Range("H31").Value = 1849.5869
$ws.Range("I31").Value = 900.5294
$ws.Range("J31").Value = 4538.5835
$ws.Range("K31").Value = 900.5294
$ws.Range("L31").Value = 4538.5835
$ws.Range("M31").Value = -605.5294
$ws.Range("N31").Value = -5128.5835
$ws.Range("H34").Value = 1849.5869
$ws.Range("I34").Value = 900.5294
$ws.Range("J34").Value = 4538.5835
$ws.Range("K34").Value = 900.5294
$ws.Range("L34").Value = 4538.5835
$ws.Range("M34").Value = -698.5294
$ws.Range("N34").Value = -4942.5835
$ws.Range("H58").Value = 589900.75
$ws.Range("I58").Value = 668236.25
$ws.Range("J58").Value = 2384.5
$ws.Range("K58").Value = 668236.25
$ws.Range("L58").Value = 2384.5
$ws.Range("M58").Value = -668033.25
$ws.Range("N58").Value = -2790.5
$ws.Range("H105").Value = 913.2273
$ws.Range("I105").Value = 804.55
$ws.Range("K105").Value = 804.55
$ws.Range("M105").Value = 942.45
$ws.Range("H107").Value = 996.2
$ws.Range("I107").Value = 981
$ws.Range("K107").Value = 981
$ws.Range("M107").Value = 939
$ws.Range("H136").Value = 589900.75
$ws.Range("I136").Value = 668236.25
$ws.Range("J136").Value = 2384.5
$ws.Range("K136").Value = 2004708.75
$ws.Range("L136").Value = 7153.5
$ws.Range("M136").Value = -2002158.75
$ws.Range("N136").Value = -12253.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 43645476
$ws.Range("I4").Value = 1167020.5
$ws.Range("J4").Value = 152875790
$ws.Range("K4").Value = 3501061.5
$ws.Range("L4").Value = 458627370
$ws.Range("M4").Value = -3500949.5
$ws.Range("N4").Value = -458627594
$ws.Range("H61").Value = 852
$ws.Range("J61").Value = 1500
$ws.Range("L61").Value = 4500
$ws.Range("N61").Value = -4930
$ws.Range("H121").Value = 1145.6
$ws.Range("J121").Value = 832.6667
$ws.Range("L121").Value = 2498.0001
$ws.Range("N121").Value = -5118.0001
$ws.Range("H139").Value = 1893
$ws.Range("I139").Value = 1800.75
$ws.Range("K139").Value = 5402.25
$ws.Range("M139").Value = -262.25
$ws.Range("H141").Value = 11449.869
$ws.Range("I141").Value = 5668.9287
$ws.Range("J141").Value = 20442.445
$ws.Range("K141").Value = 17006.7861
$ws.Range("L141").Value = 61327.335
$ws.Range("M141").Value = -11826.7861
$ws.Range("N141").Value = -71687.33499999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 20000
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20596
$ws.Range("H97").Value = 595.44446
$ws.Range("I97").Value = 683.3077
$ws.Range("J97").Value = 367
$ws.Range("K97").Value = 683.3077
$ws.Range("L97").Value = 367
$ws.Range("M97").Value = -187.3077
$ws.Range("N97").Value = -1359

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5000
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -4888
$ws.Range("N2").Value = -5224
$ws.Range("H22").Value = 1458.1
$ws.Range("I22").Value = 824.5
$ws.Range("J22").Value = 1616.5
$ws.Range("K22").Value = 824.5
$ws.Range("L22").Value = 1616.5
$ws.Range("M22").Value = -529.5
$ws.Range("N22").Value = -2206.5
$ws.Range("H27").Value = 1458.1
$ws.Range("I27").Value = 824.5
$ws.Range("J27").Value = 1616.5
$ws.Range("K27").Value = 824.5
$ws.Range("L27").Value = 1616.5
$ws.Range("M27").Value = -717.5
$ws.Range("N27").Value = -1830.5
$ws.Range("H40").Value = 4842.0557
$ws.Range("I40").Value = 3699.6667
$ws.Range("K40").Value = 3699.6667
$ws.Range("M40").Value = -3563.6667
$ws.Range("H82").Value = 2842869.2
$ws.Range("I82").Value = 5210250.5
$ws.Range("J82").Value = 2011.8
$ws.Range("K82").Value = 5210250.5
$ws.Range("L82").Value = 2011.8
$ws.Range("M82").Value = -5209889.5
$ws.Range("N82").Value = -2733.8
$ws.Range("H85").Value = 2842869.2
$ws.Range("I85").Value = 5210250.5
$ws.Range("J85").Value = 2011.8
$ws.Range("K85").Value = 5210250.5
$ws.Range("L85").Value = 2011.8
$ws.Range("M85").Value = -5209002.5
$ws.Range("N85").Value = -4507.8
$ws.Range("H93").Value = 1129.6
$ws.Range("I93").Value = 941.1667
$ws.Range("K93").Value = 941.1667
$ws.Range("M93").Value = 306.8333
$ws.Range("H100").Value = 4995.6665
$ws.Range("J100").Value = 4998.5
$ws.Range("L100").Value = 4998.5
$ws.Range("N100").Value = -6080.5
$ws.Range("H122").Value = 87915790
$ws.Range("I122").Value = 111114024
$ws.Range("K122").Value = 333342072
$ws.Range("M122").Value = -333339622
$ws.Range("H132").Value = 4455.68
$ws.Range("I132").Value = 3477.1943
$ws.Range("K132").Value = 10431.5829
$ws.Range("M132").Value = -7901.582900000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 75626
$ws.Range("J27").Value = 75626
$ws.Range("L27").Value = 75626
$ws.Range("N27").Value = -75764
$ws.Range("H100").Value = 1669621.1
$ws.Range("I100").Value = 4002218
$ws.Range("K100").Value = 8004436
$ws.Range("M100").Value = -8003895
$ws.Range("H107").Value = 2073.6453
$ws.Range("I107").Value = 2327.64
$ws.Range("J107").Value = 1015.3333
$ws.Range("K107").Value = 6982.92
$ws.Range("L107").Value = 3045.9999
$ws.Range("M107").Value = -5062.92
$ws.Range("N107").Value = -6885.9999
$ws.Range("H132").Value = 15155304
$ws.Range("I132").Value = 1634.9
$ws.Range("K132").Value = 4904.700000000001
$ws.Range("M132").Value = -2374.700000000001

Write-Host "Applied 256 cell updates."